# chore: update Sheets via scheduled runner
# Refresh market-price / leve-profit columns (H:N) with latest pulled data
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR leve-profit tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6761.36
$ws.Range("I28").Value = 810.5714
$ws.Range("J28").Value = 38003
$ws.Range("K28").Value = 810.5714
$ws.Range("L28").Value = 38003
$ws.Range("M28").Value = -325.5714
$ws.Range("N28").Value = -38973
$ws.Range("H40").Value = 918.8570999999999
$ws.Range("I40").Value = 663.5
$ws.Range("J40").Value = 1110.375
$ws.Range("K40").Value = 663.5
$ws.Range("L40").Value = 1110.375
$ws.Range("M40").Value = -488.5
$ws.Range("N40").Value = -1460.375
$ws.Range("H53").Value = 188.15384
$ws.Range("I53").Value = 260.25
$ws.Range("J53").Value = 156.11111
$ws.Range("K53").Value = 260.25
$ws.Range("L53").Value = 156.11111
$ws.Range("M53").Value = 376.75
$ws.Range("N53").Value = -1430.11111
$ws.Range("H129").Value = 1339.3948
$ws.Range("J129").Value = 1362.081
$ws.Range("L129").Value = 4086.242999999999
$ws.Range("N129").Value = -14086.243
$ws.Range("H132").Value = 2229.22
$ws.Range("I132").Value = 2189.7659
$ws.Range("J132").Value = 2847.3333
$ws.Range("K132").Value = 6569.297699999999
$ws.Range("L132").Value = 8541.999899999999
$ws.Range("M132").Value = -4039.297699999999
$ws.Range("N132").Value = -13601.9999
$ws.Range("H138").Value = 4550.607
$ws.Range("I138").Value = 5997.3335
$ws.Range("J138").Value = 4377
$ws.Range("K138").Value = 17992.0005
$ws.Range("L138").Value = 13131
$ws.Range("M138").Value = -12852.0005
$ws.Range("N138").Value = -23411

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20429608
$ws.Range("I32").Value = 24409948
$ws.Range("J32").Value = 30375
$ws.Range("K32").Value = 24409948
$ws.Range("L32").Value = 30375
$ws.Range("M32").Value = -24409661
$ws.Range("N32").Value = -30949
$ws.Range("H74").Value = 12821859
$ws.Range("I74").Value = 1127.4722
$ws.Range("J74").Value = 166670640
$ws.Range("K74").Value = 1127.4722
$ws.Range("L74").Value = 166670640
$ws.Range("M74").Value = -253.4721999999999
$ws.Range("N74").Value = -166672388
$ws.Range("H77").Value = 12821859
$ws.Range("I77").Value = 1127.4722
$ws.Range("J77").Value = 166670640
$ws.Range("K77").Value = 5637.361
$ws.Range("L77").Value = 833353200
$ws.Range("M77").Value = -1269.361
$ws.Range("N77").Value = -833361936
$ws.Range("H132").Value = 1792033
$ws.Range("I132").Value = 2979.4243
$ws.Range("J132").Value = 7695910
$ws.Range("K132").Value = 8938.2729
$ws.Range("L132").Value = 23087730
$ws.Range("M132").Value = -6408.2729
$ws.Range("N132").Value = -23092790

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3333634.8
$ws.Range("I7").Value = 5000150
$ws.Range("J7").Value = 604
$ws.Range("K7").Value = 5000150
$ws.Range("L7").Value = 604
$ws.Range("M7").Value = -5000037
$ws.Range("N7").Value = -830
$ws.Range("H134").Value = 4131.931
$ws.Range("I134").Value = 4180.25
$ws.Range("J134").Value = 3900
$ws.Range("K134").Value = 12540.75
$ws.Range("L134").Value = 11700
$ws.Range("M134").Value = -10005.75
$ws.Range("N134").Value = -16770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6980.1196
$ws.Range("I31").Value = 995.38464
$ws.Range("J31").Value = 9337.742
$ws.Range("K31").Value = 995.38464
$ws.Range("L31").Value = 9337.742
$ws.Range("M31").Value = -700.38464
$ws.Range("N31").Value = -9927.742
$ws.Range("H34").Value = 6980.1196
$ws.Range("I34").Value = 995.38464
$ws.Range("J34").Value = 9337.742
$ws.Range("K34").Value = 995.38464
$ws.Range("L34").Value = 9337.742
$ws.Range("M34").Value = -793.38464
$ws.Range("N34").Value = -9741.742
$ws.Range("H58").Value = 2521.6365
$ws.Range("I58").Value = 2304.2222
$ws.Range("K58").Value = 2304.2222
$ws.Range("M58").Value = -2101.2222
$ws.Range("H132").Value = 55557940
$ws.Range("J132").Value = 33335834
$ws.Range("L132").Value = 100007502
$ws.Range("N132").Value = -100012562
$ws.Range("H136").Value = 2521.6365
$ws.Range("I136").Value = 2304.2222
$ws.Range("K136").Value = 6912.6666
$ws.Range("M136").Value = -4362.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1142.127
$ws.Range("I5").Value = 871.9623
$ws.Range("J5").Value = 2574
$ws.Range("K5").Value = 2615.8869
$ws.Range("L5").Value = 7722
$ws.Range("M5").Value = -2503.8869
$ws.Range("N5").Value = -7946
$ws.Range("H50").Value = 322.92307
$ws.Range("I50").Value = 226.54546
$ws.Range("J50").Value = 853
$ws.Range("K50").Value = 679.6363799999999
$ws.Range("L50").Value = 2559
$ws.Range("M50").Value = -198.6363799999999
$ws.Range("N50").Value = -3521
$ws.Range("H53").Value = 322.92307
$ws.Range("I53").Value = 226.54546
$ws.Range("J53").Value = 853
$ws.Range("K53").Value = 679.6363799999999
$ws.Range("L53").Value = 2559
$ws.Range("M53").Value = -198.6363799999999
$ws.Range("N53").Value = -3521
$ws.Range("H113").Value = 658.1719000000001
$ws.Range("I113").Value = 643.6585
$ws.Range("J113").Value = 684.04346
$ws.Range("K113").Value = 1930.9755
$ws.Range("L113").Value = 2052.13038
$ws.Range("M113").Value = 239.0245
$ws.Range("N113").Value = -6392.130380000001
$ws.Range("H129").Value = 1685287.9
$ws.Range("J129").Value = 2333256.8
$ws.Range("L129").Value = 6999770.399999999
$ws.Range("N129").Value = -7009770.399999999
$ws.Range("H131").Value = 2933.5078
$ws.Range("I131").Value = 857.1429000000001
$ws.Range("J131").Value = 3184.1035
$ws.Range("K131").Value = 2571.4287
$ws.Range("L131").Value = 9552.3105
$ws.Range("M131").Value = 2468.5713
$ws.Range("N131").Value = -19632.3105
$ws.Range("H135").Value = 1142.127
$ws.Range("I135").Value = 871.9623
$ws.Range("J135").Value = 2574
$ws.Range("K135").Value = 7847.6607
$ws.Range("L135").Value = 23166
$ws.Range("M135").Value = -5312.6607
$ws.Range("N135").Value = -28236
$ws.Range("H136").Value = 3835.7917
$ws.Range("I136").Value = 1454.9166
$ws.Range("J136").Value = 6216.6665
$ws.Range("K136").Value = 4364.7498
$ws.Range("L136").Value = 18649.9995
$ws.Range("M136").Value = 735.2502000000004
$ws.Range("N136").Value = -28849.9995
$ws.Range("H137").Value = 10138255
$ws.Range("J137").Value = 14540457
$ws.Range("L137").Value = 43621371
$ws.Range("N137").Value = -43631571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9248
$ws.Range("J2").Value = 9248
$ws.Range("L2").Value = 9248
$ws.Range("N2").Value = -9472
$ws.Range("H46").Value = 519.3333
$ws.Range("I46").Value = 476
$ws.Range("J46").Value = 550.2857
$ws.Range("K46").Value = 476
$ws.Range("L46").Value = 550.2857
$ws.Range("M46").Value = -288
$ws.Range("N46").Value = -926.2857
$ws.Range("H55").Value = 507
$ws.Range("J55").Value = 876.4
$ws.Range("L55").Value = 876.4
$ws.Range("N55").Value = -1222.4
$ws.Range("H132").Value = 3518.875
$ws.Range("I132").Value = 2935.3572
$ws.Range("J132").Value = 3972.7222
$ws.Range("K132").Value = 8806.071599999999
$ws.Range("L132").Value = 11918.1666
$ws.Range("M132").Value = -6276.071599999999
$ws.Range("N132").Value = -16978.1666
$ws.Range("H136").Value = 3624870.2
$ws.Range("I136").Value = 1164.8387
$ws.Range("J136").Value = 11113862
$ws.Range("K136").Value = 3494.5161
$ws.Range("L136").Value = 33341586
$ws.Range("M136").Value = -944.5160999999998
$ws.Range("N136").Value = -33346686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 45151.43
$ws.Range("J22").Value = 45151.43
$ws.Range("L22").Value = 45151.43
$ws.Range("N22").Value = -45737.43
$ws.Range("H28").Value = 90018.75
$ws.Range("I28").Value = 60017
$ws.Range("J28").Value = 100019.336
$ws.Range("K28").Value = 60017
$ws.Range("L28").Value = 100019.336
$ws.Range("M28").Value = -59669
$ws.Range("N28").Value = -100715.336
$ws.Range("H122").Value = 2841.05
$ws.Range("I122").Value = 2415.4285
$ws.Range("K122").Value = 7246.2855
$ws.Range("M122").Value = -4796.2855
$ws.Range("H126").Value = 2175.1333
$ws.Range("I126").Value = 1969.75
$ws.Range("J126").Value = 2996.6667
$ws.Range("K126").Value = 5909.25
$ws.Range("L126").Value = 8990.000100000001
$ws.Range("M126").Value = -3439.25
$ws.Range("N126").Value = -13930.0001
